$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (shifts the existing Student Id / Pathway Date /
# Order / Screen / Message / Url columns one place to the right) and give it
# a header of "Student Name", matching the style used by the other header
# cells (bold Verdana).
$ws.Range("A1").EntireColumn.Insert() | Out-Null

$ws.Range("A1").Value = "Student Name"
$ws.Range("A1").Font.Bold = $true

# New column gets its own explicit width.
$ws.Columns.Item(1).ColumnWidth = 11.25

# Leave the selection where the editor ended up after adding the column.
$ws.Range("B3").Select() | Out-Null
